# Commit: "update the number of significant digits"
# Round a handful of measured values (vessel diameter / density) down to
# fewer significant digits on the "Vessel size (tumor)" and
# "Vessel density (tumor)" sheets.

$wb = $excel.ActiveWorkbook

# --- Vessel size (tumor) ---------------------------------------------
$wsSize = $wb.Worksheets.Item("Vessel size (tumor)")
$wsSize.Range("C3").Value = 15        # was 14.76
$wsSize.Range("B4").Value = 86.7      # was 86.66
$wsSize.Range("C4").Value = 2.76      # was 2.755
$wsSize.Range("B6").Value = 135       # was 134.8
$wsSize.Range("B8").Value = 113       # was 113.1

# leave the selection where the edits ended, matching the saved view
$wsSize.Range("B8").Select()

# --- Vessel density (tumor) -------------------------------------------
$wsDensity = $wb.Worksheets.Item("Vessel density (tumor)")
$wsDensity.Range("B2").Value = 19.5    # was 19.45
$wsDensity.Range("B3").Value = 19.1    # was 19.13
$wsDensity.Range("B4").Value = 17.6    # was 17.64
$wsDensity.Range("B11").Value = 292    # was 292.45
$wsDensity.Range("C11").Value = 28.6   # was 28.64
$wsDensity.Range("B12").Value = 212    # was 211.93

$wsDensity.Activate()
$wsDensity.Range("B12").Select()
